# Auto-generated edit script: updates computed price/profit columns (H-N)
# across multiple worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19: Unbreak My Heart - Roof Tile
$ws.Range("H19").Value = 4446
$ws.Range("I19").Value = 4425.5
$ws.Range("K19").Value = 4425.5
$ws.Range("M19").Value = -4250.5

# ALC row 33: Glazed and Confused - Clear Glass Lens
$ws.Range("H33").Value = 338.53333
$ws.Range("I33").Value = 148.42857
$ws.Range("K33").Value = 148.42857
$ws.Range("M33").Value = 80.57142999999999

# ALC row 40: Stuck in the Moment - Horn Glue
$ws.Range("H40").Value = 7027.5
$ws.Range("J40").Value = 8142.857
$ws.Range("L40").Value = 8142.857
$ws.Range("N40").Value = -8492.857

# ALC row 43: Growing Is Knowing - Growth Formula Gamma
$ws.Range("H43").Value = 4378.9
$ws.Range("J43").Value = 5597.75
$ws.Range("L43").Value = 5597.75
$ws.Range("N43").Value = -5735.75

# ALC row 132: Fast-forwarding Flora - Growth Formula Lambda
$ws.Range("H132").Value = 2193.1667
$ws.Range("I132").Value = 1940.7391
$ws.Range("K132").Value = 5822.2173
$ws.Range("M132").Value = -3292.2173

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth - Bronze Rivets
$ws.Range("H5").Value = 208.72728
$ws.Range("I5").Value = 135.22223
$ws.Range("K5").Value = 135.22223
$ws.Range("M5").Value = -23.22223

# ARM row 45: Hollow Hallmarks - Mythril Ingot
$ws.Range("H45").Value = 340916.66
$ws.Range("I45").Value = 1003750
$ws.Range("J45").Value = 9500
$ws.Range("K45").Value = 1003750
$ws.Range("L45").Value = 9500
$ws.Range("M45").Value = -1003373
$ws.Range("N45").Value = -10254

# ARM row 46: Get Me the Usual - Heavy Steel Flanchard
$ws.Range("H46").Value = 7259
$ws.Range("J46").Value = 7410
$ws.Range("L46").Value = 7410
$ws.Range("N46").Value = -8048

# ARM row 61: Dealing with the Tough Stuff - Cobalt Ingot
$ws.Range("H61").Value = 7024.074
$ws.Range("I61").Value = 5394.591
$ws.Range("K61").Value = 5394.591
$ws.Range("M61").Value = -5182.591

# ARM row 102: Smells of Rich Tama-hagane - Tama-hagane Ingot
$ws.Range("H102").Value = 3384.276
$ws.Range("I102").Value = 2447
$ws.Range("J102").Value = 9242.25
$ws.Range("K102").Value = 2447
$ws.Range("L102").Value = 9242.25
$ws.Range("M102").Value = -825
$ws.Range("N102").Value = -12486.25

# ARM row 110: Scheduled Maintenance - Deepgold Ingot
$ws.Range("H110").Value = 6064.421
$ws.Range("I110").Value = 4435.5835
$ws.Range("J110").Value = 8856.714
$ws.Range("K110").Value = 4435.5835
$ws.Range("L110").Value = 8856.714
$ws.Range("M110").Value = -2390.5835
$ws.Range("N110").Value = -12946.714

# ARM row 132: Don't Bore Me, Ore Me - Mountain Chromite Ingot
$ws.Range("H132").Value = 4414.1143
$ws.Range("J132").Value = 8849
$ws.Range("L132").Value = 26547
$ws.Range("N132").Value = -31607

# ARM row 136: Metal with Mettle - Cobalt Tungsten Ingot
$ws.Range("H136").Value = 7024.074
$ws.Range("I136").Value = 5394.591
$ws.Range("K136").Value = 16183.773
$ws.Range("M136").Value = -13633.773

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences - Bronze Rivets
$ws.Range("H4").Value = 208.72728
$ws.Range("I4").Value = 135.22223
$ws.Range("K4").Value = 135.22223
$ws.Range("M4").Value = -20.22223

# BSM row 22: Riveting Run - Iron Rivets
$ws.Range("H22").Value = 1563.7333
$ws.Range("I22").Value = 2221.4285
$ws.Range("K22").Value = 2221.4285
$ws.Range("M22").Value = -2048.4285

# BSM row 64: With Bearings Straight - Mythrite Nugget
$ws.Range("H64").Value = 709.8570999999999
$ws.Range("I64").Value = 664.6667
$ws.Range("J64").Value = 743.75
$ws.Range("K64").Value = 664.6667
$ws.Range("L64").Value = 743.75
$ws.Range("M64").Value = -439.6667
$ws.Range("N64").Value = -1193.75

# BSM row 67: Bearing the Brunt (L) - Mythrite Nugget
$ws.Range("H67").Value = 709.8570999999999
$ws.Range("I67").Value = 664.6667
$ws.Range("J67").Value = 743.75
$ws.Range("K67").Value = 664.6667
$ws.Range("L67").Value = 743.75
$ws.Range("M67").Value = 115.3333
$ws.Range("N67").Value = -2303.75

# BSM row 86: Through Thick and Thin - Adamantite Nugget
$ws.Range("H86").Value = 1980.4166
$ws.Range("I86").Value = 1751.3636
$ws.Range("K86").Value = 1751.3636
$ws.Range("M86").Value = -628.3635999999999

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) - Adamantite Nugget
$ws.Range("H89").Value = 1980.4166
$ws.Range("I89").Value = 1751.3636
$ws.Range("K89").Value = 8756.817999999999
$ws.Range("M89").Value = -3140.817999999999

# BSM row 94: High Steal - High Steel Nugget
$ws.Range("H94").Value = 4264.933
$ws.Range("I94").Value = 3497.9167
$ws.Range("K94").Value = 3497.9167
$ws.Range("M94").Value = -3046.9167

# BSM row 134: Ruthenium Supremium - Ruthenium Ingot
$ws.Range("H134").Value = 3793.3928
$ws.Range("I134").Value = 3793.3928
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11380.1784
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8845.178400000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found - Walnut Lumber
$ws.Range("H31").Value = 4047.276
$ws.Range("I31").Value = 2816.3333
$ws.Range("K31").Value = 2816.3333
$ws.Range("M31").Value = -2521.3333

# CRP row 34: Armoires of the Rich and Famous - Walnut Lumber
$ws.Range("H34").Value = 4047.276
$ws.Range("I34").Value = 2816.3333
$ws.Range("K34").Value = 2816.3333
$ws.Range("M34").Value = -2614.3333

# CRP row 41: The Lone Bowman - Oak Longbow
$ws.Range("H41").Value = 20940.555
$ws.Range("I41").Value = 20849.5
$ws.Range("K41").Value = 20849.5
$ws.Range("M41").Value = -20421.5

# CRP row 47: Grippy When Wet - Mythril Cavalry Bow
$ws.Range("H47").Value = 30000.5
$ws.Range("I47").Value = 30000
$ws.Range("K47").Value = 30000
$ws.Range("M47").Value = -29434

# CRP row 59: Bow Down to Magic - Crab Bow
$ws.Range("H59").Value = 40921.75
$ws.Range("J59").Value = 40979
$ws.Range("L59").Value = 40979
$ws.Range("N59").Value = -43269

# CRP row 74: License to Heal - Dark Chestnut Rod
$ws.Range("H74").Value = 35021.43
$ws.Range("J74").Value = 35021.43
$ws.Range("L74").Value = 35021.43
$ws.Range("N74").Value = -36769.43

# CRP row 77: Purified Polyrhythm (L) - Dark Chestnut Rod
$ws.Range("H77").Value = 35021.43
$ws.Range("J77").Value = 35021.43
$ws.Range("L77").Value = 105064.29
$ws.Range("N77").Value = -113800.29

# CRP row 99: O Pine - Pine Lumber
$ws.Range("H99").Value = 2900.3333
$ws.Range("I99").Value = 2900.3333
$ws.Range("K99").Value = 2900.3333
$ws.Range("M99").Value = -1402.3333

# CRP row 105: Zelkova, My Love - Zelkova Lumber
$ws.Range("H105").Value = 1999.2858
$ws.Range("I105").Value = 2046.5385
$ws.Range("K105").Value = 2046.5385
$ws.Range("M105").Value = -299.5385000000001

# CRP row 126: A Better Conductor - Red Pine Lumber
$ws.Range("H126").Value = 2900.3333
$ws.Range("I126").Value = 2900.3333
$ws.Range("K126").Value = 8700.999899999999
$ws.Range("M126").Value = -6230.999899999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 23: Sweet Smell of Success - Lavender Oil
$ws.Range("H23").Value = 11180.223
$ws.Range("I23").Value = 60.42857
$ws.Range("K23").Value = 181.28571
$ws.Range("M23").Value = 53.71429000000001

# CUL row 33: Cooking with Gas - Chicken Stock
$ws.Range("H33").Value = 541.1667
$ws.Range("I33").Value = 43.75
$ws.Range("J33").Value = 1536
$ws.Range("K33").Value = 262.5
$ws.Range("L33").Value = 9216
$ws.Range("M33").Value = 20.5
$ws.Range("N33").Value = -9782

# CUL row 80: Saucy for a Suitor - Hollandaise Sauce
$ws.Range("H80").Value = 5999.5
$ws.Range("J80").Value = 5999
$ws.Range("L80").Value = 17997
$ws.Range("N80").Value = -19869

# CUL row 83: Saved by the Sauce (L) - Hollandaise Sauce
$ws.Range("H83").Value = 5999.5
$ws.Range("J83").Value = 5999
$ws.Range("L83").Value = 53991
$ws.Range("N83").Value = -63351

# CUL row 131: The Mountain Steeped - Tsai tou Vounou
$ws.Range("H131").Value = 38464230
$ws.Range("J131").Value = 2919.9167
$ws.Range("L131").Value = 8759.750100000001
$ws.Range("N131").Value = -18839.7501

$ws = $wb.Worksheets.Item("GSM")
# GSM row 10: Let's Talk about Hex - Bone Necklace
$ws.Range("H10").Value = 312.5
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("N10").Value = -838

$ws = $wb.Worksheets.Item("LTW")
# LTW row 2: Red in the Head - Leather Calot
$ws.Range("H2").Value = 10000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224

# LTW row 16: Saddle Sore - Hard Leather
$ws.Range("H16").Value = 1678.0667
$ws.Range("I16").Value = 1404.5454
$ws.Range("J16").Value = 2430.25
$ws.Range("K16").Value = 1404.5454
$ws.Range("L16").Value = 2430.25
$ws.Range("M16").Value = -1234.5454
$ws.Range("N16").Value = -2770.25

# LTW row 46: Supply Side Logic - Boar Leather
$ws.Range("H46").Value = 15209.2
$ws.Range("I46").Value = 3124.75
$ws.Range("K46").Value = 3124.75
$ws.Range("M46").Value = -2936.75

# LTW row 55: It's Not a Job, It's a Calling - Peiste Leather
$ws.Range("H55").Value = 1165
$ws.Range("I55").Value = 1183.3334
$ws.Range("J55").Value = 1137.5
$ws.Range("K55").Value = 1183.3334
$ws.Range("L55").Value = 1137.5
$ws.Range("M55").Value = -1010.3334
$ws.Range("N55").Value = -1483.5

# LTW row 122: Hell on Leather - Gaja Leather
$ws.Range("H122").Value = 3393.9443
$ws.Range("I122").Value = 3352.4
$ws.Range("K122").Value = 10057.2
$ws.Range("M122").Value = -7607.200000000001

# LTW row 132: Tenets of Tanning - Silver Lobo Leather
$ws.Range("H132").Value = 9790.5
$ws.Range("J132").Value = 9497.75
$ws.Range("L132").Value = 28493.25
$ws.Range("N132").Value = -33553.25

$ws = $wb.Worksheets.Item("WVR")
# WVR row 107: Flax Wax - Bright Linen Yarn
$ws.Range("H107").Value = 585.3333
$ws.Range("I107").Value = 532.5
$ws.Range("K107").Value = 1597.5
$ws.Range("M107").Value = 322.5

# WVR row 132: Comfy Cabins - Snow Cotton Cloth
$ws.Range("H132").Value = 3615.68
$ws.Range("I132").Value = 3126.85
$ws.Range("K132").Value = 9380.549999999999
$ws.Range("M132").Value = -6850.549999999999

# WVR row 136: Weaving the Envelope - Sarcenet Cloth
$ws.Range("H136").Value = 3304.6538
$ws.Range("J136").Value = 5564.2
$ws.Range("L136").Value = 16692.6
$ws.Range("N136").Value = -21792.6
